$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (new Price text, new Volume(1h) text) for the rows that changed.
$updates = @{
  2 = @("27.343.36", "  +1.15%  ")
  3 = @("1.856.75", "  +1.65%  ")
  4 = @("'1.002", "  -0.53%  ")
  5 = @("'314.44", "  +1.46%  ")
  6 = @("'1.001", "  -0.50%  ")
  7 = @("'0.4625", "  -0.55%  ")
  8 = @("'0.3713", "  +1.51%  ")
  9 = @("'0.07342", "  +1.60%  ")
  10 = @("'0.8824", "  +2.59%  ")
  11 = @("'0.07881", "  +1.74%  ")
  12 = @("'19.86", "  -0.09%  ")
  13 = @("1.916.18", "  +7.17%  ")
  14 = @("'5.390", "  +1.15%  ")
  15 = @("'6.554", "  +0.88%  ")
  16 = @("'92.04", "  +0.30%  ")
  17 = @("'1.003", "  -0.46%  ")
  18 = @($null, "  +2.60%  ")
  19 = @("'1.002", "  -0.37%  ")
  20 = @($null, "  +2.42%  ")
  21 = @("27.371.92", "  +2.07%  ")
  22 = @("'5.122", "  -0.57%  ")
  23 = @("'10.52", "  -0.03%  ")
  24 = @("2.115.93", "  +0.53%  ")
  25 = @("'152.65", "  +0.61%  ")
  26 = @("'1.880", "  +2.33%  ")
  27 = @("'18.39", "  +1.10%  ")
  28 = @("'2.077", "  +0.95%  ")
  29 = @("'5.131", "  +0.55%  ")
  30 = @("'116.27", "  +0.75%  ")
  31 = @("'0.08886", "  +0.89%  ")
  32 = @($null, "  +5.23%  ")
  33 = @("'3.025", "  +2.42%  ")
  34 = @("'1.166", "  +3.14%  ")
  35 = @("'4.489", "  +1.55%  ")
  36 = @("'2.611", "  +8.58%  ")
  37 = @($null, "  +0.30%  ")
  38 = @($null, "  +1.65%  ")
  39 = @("'2.975", "  +1.47%  ")
  40 = @("'0.05227", "  -0.20%  ")
  41 = @("'7.078", "  -0.95%  ")
  42 = @("'0.5161", "  +0.04%  ")
  43 = @("'0.1644", "  +1.06%  ")
  44 = @("'8.365", "  +2.26%  ")
  45 = @($null, "  +1.07%  ")
  46 = @("'10.27", "  +1.44%  ")
  47 = @("'1.001", "  -0.56%  ")
  48 = @("'103.34", "  +0.59%  ")
  49 = @("'1.654", "  +2.31%  ")
  50 = @("'0.06228", "  -0.22%  ")
  51 = @("'65.62", "  +2.61%  ")
}

foreach ($r in $updates.Keys) {
  $pair = $updates[$r]
  if ($pair[0] -ne $null) {
    $ws.Cells.Item($r, 4).Value = $pair[0]
  }
  $ws.Cells.Item($r, 5).Value = $pair[1]
}
